# Auto-generated edit script: update Leve profit values across sheets
# per scheduled runner refresh of market board prices (Aegis_Profits).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 524
$ws.Range("I92").Value = 436.33334
$ws.Range("K92").Value = 436.33334
$ws.Range("M92").Value = 811.66666
# Row 103
$ws.Range("H103").Value = 1344.6666
$ws.Range("J103").Value = 760.9091
$ws.Range("L103").Value = 2282.7273
$ws.Range("N103").Value = -3454.7273

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25858.027
$ws.Range("I32").Value = 5902.746
$ws.Range("J32").Value = 151576.3
$ws.Range("K32").Value = 5902.746
$ws.Range("L32").Value = 151576.3
$ws.Range("M32").Value = -5615.746
$ws.Range("N32").Value = -152150.3
# Row 37
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4727
# Row 45
$ws.Range("H45").Value = 85612.664
$ws.Range("I45").Value = 112432.445
$ws.Range("J45").Value = 5153.3335
$ws.Range("K45").Value = 112432.445
$ws.Range("L45").Value = 5153.3335
$ws.Range("M45").Value = -112055.445
$ws.Range("N45").Value = -5907.3335
# Row 122
$ws.Range("H122").Value = 3624.2856
$ws.Range("I122").Value = 3624.2856
$ws.Range("K122").Value = 10872.8568
$ws.Range("M122").Value = -8422.856800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 47394
$ws.Range("J62").Value = 47394
$ws.Range("L62").Value = 47394
$ws.Range("N62").Value = -48766
# Row 65
$ws.Range("H65").Value = 47394
$ws.Range("J65").Value = 47394
$ws.Range("L65").Value = 142182
$ws.Range("N65").Value = -149046

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1210.5883
$ws.Range("I6").Value = 1049.1666
$ws.Range("J6").Value = 1598
$ws.Range("K6").Value = 1049.1666
$ws.Range("L6").Value = 1598
$ws.Range("M6").Value = -936.1666
$ws.Range("N6").Value = -1824
# Row 7
$ws.Range("H7").Value = 97.666664
$ws.Range("I7").Value = 63.11111
$ws.Range("J7").Value = 201.33333
$ws.Range("K7").Value = 63.11111
$ws.Range("L7").Value = 201.33333
$ws.Range("M7").Value = 49.88889
$ws.Range("N7").Value = -427.33333
# Row 31
$ws.Range("H31").Value = 26212.05
$ws.Range("I31").Value = 1258.7037
$ws.Range("J31").Value = 47945.613
$ws.Range("K31").Value = 1258.7037
$ws.Range("L31").Value = 47945.613
$ws.Range("M31").Value = -963.7037
$ws.Range("N31").Value = -48535.613
# Row 34
$ws.Range("H34").Value = 26212.05
$ws.Range("I34").Value = 1258.7037
$ws.Range("J34").Value = 47945.613
$ws.Range("K34").Value = 1258.7037
$ws.Range("L34").Value = 47945.613
$ws.Range("M34").Value = -1056.7037
$ws.Range("N34").Value = -48349.613
# Row 41
$ws.Range("H41").Value = 12112.714
$ws.Range("I41").Value = 5050
$ws.Range("J41").Value = 13289.833
$ws.Range("K41").Value = 5050
$ws.Range("L41").Value = 13289.833
$ws.Range("M41").Value = -4622
$ws.Range("N41").Value = -14145.833
# Row 50
$ws.Range("H50").Value = 14546.667
$ws.Range("J50").Value = 14546.667
$ws.Range("L50").Value = 14546.667
$ws.Range("N50").Value = -15796.667
# Row 51
$ws.Range("H51").Value = 9384.385
$ws.Range("J51").Value = 9384.385
$ws.Range("L51").Value = 9384.385
$ws.Range("N51").Value = -10856.385
# Row 59
$ws.Range("H59").Value = 33653.332
$ws.Range("J59").Value = 33653.332
$ws.Range("L59").Value = 33653.332
$ws.Range("N59").Value = -35943.332
# Row 60
$ws.Range("H60").Value = 16740
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
# Row 61
$ws.Range("H61").Value = 9384.385
$ws.Range("J61").Value = 9384.385
$ws.Range("L61").Value = 9384.385
$ws.Range("N61").Value = -10080.385
# Row 68
$ws.Range("H68").Value = 18007.908
$ws.Range("J68").Value = 18007.908
$ws.Range("L68").Value = 18007.908
$ws.Range("N68").Value = -19505.908
# Row 71
$ws.Range("H71").Value = 18007.908
$ws.Range("J71").Value = 18007.908
$ws.Range("L71").Value = 54023.724
$ws.Range("N71").Value = -61511.724
# Row 74
$ws.Range("H74").Value = 29427
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 29427
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 29427
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -31175
# Row 77
$ws.Range("H77").Value = 29427
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 29427
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 88281
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -97017
# Row 86
$ws.Range("H86").Value = 2619
$ws.Range("I86").Value = 2466.5
$ws.Range("J86").Value = 2695.25
$ws.Range("K86").Value = 2466.5
$ws.Range("L86").Value = 2695.25
$ws.Range("M86").Value = -1343.5
$ws.Range("N86").Value = -4941.25
# Row 89
$ws.Range("H89").Value = 2619
$ws.Range("I89").Value = 2466.5
$ws.Range("J89").Value = 2695.25
$ws.Range("K89").Value = 12332.5
$ws.Range("L89").Value = 13476.25
$ws.Range("M89").Value = -6716.5
$ws.Range("N89").Value = -24708.25
# Row 135
$ws.Range("H135").Value = 46506
$ws.Range("J135").Value = 49560
$ws.Range("L135").Value = 49560
$ws.Range("N135").Value = -59700

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 326.7
$ws.Range("I40").Value = 140.55556
$ws.Range("K40").Value = 562.2222400000001
$ws.Range("M40").Value = -493.2222400000001
# Row 68
$ws.Range("H68").Value = 2453.2
$ws.Range("I68").Value = 400
$ws.Range("K68").Value = 1200
$ws.Range("M68").Value = -389
# Row 71
$ws.Range("H71").Value = 2453.2
$ws.Range("I71").Value = 400
$ws.Range("K71").Value = 3600
$ws.Range("M71").Value = 456
# Row 131
$ws.Range("H131").Value = 15205.19
$ws.Range("J131").Value = 16220.743
$ws.Range("L131").Value = 48662.229
$ws.Range("N131").Value = -58742.229

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 21900
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 21900
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 21900
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -23272
# Row 66
$ws.Range("H66").Value = 21900
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 21900
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 65700
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -72564
# Row 122
$ws.Range("H122").Value = 799.6667
$ws.Range("I122").Value = 859.6
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 2578.8
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -128.8000000000002
$ws.Range("N122").Value = -6400

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2177.25
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 1903
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 5709
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -10609
# Row 127
$ws.Range("H127").Value = 20000
$ws.Range("J127").Value = 20000
$ws.Range("L127").Value = 20000
$ws.Range("N127").Value = -29920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 643.35187
$ws.Range("I136").Value = 366.75
$ws.Range("J136").Value = 1860.4
$ws.Range("K136").Value = 1100.25
$ws.Range("L136").Value = 5581.200000000001
$ws.Range("M136").Value = 1449.75
$ws.Range("N136").Value = -10681.2
